$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F7").Value = 817
$ws1.Range("F10").Value = 2062
$ws1.Range("F16").Value = 2086
$ws1.Range("F18").Value = 9602
$ws1.Range("F19").Value = 936

# Sheet "演出" (Performances)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F12").Value = 68

# Sheet "全部类型" (All types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F12").Value = 817
$ws4.Range("F16").Value = 2062
$ws4.Range("F27").Value = 2086
$ws4.Range("F29").Value = 68
$ws4.Range("F31").Value = 936
